$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H74").Value = 8175.6665
$ws.Range("I74").Value = 6297.0586
$ws.Range("K74").Value = 6297.0586
$ws.Range("M74").Value = -5361.0586
$ws.Range("H77").Value = 8175.6665
$ws.Range("I77").Value = 6297.0586
$ws.Range("K77").Value = 31485.293
$ws.Range("M77").Value = -26805.293
$ws.Range("H98").Value = 938.5333000000001
$ws.Range("I98").Value = 950.0833
$ws.Range("J98").Value = 892.3333
$ws.Range("K98").Value = 950.0833
$ws.Range("L98").Value = 892.3333
$ws.Range("M98").Value = 547.9167
$ws.Range("N98").Value = -3888.3333
$ws.Range("H116").Value = 3676.6667
$ws.Range("I116").Value = 3398.3333
$ws.Range("J116").Value = 3815.8333
$ws.Range("K116").Value = 3398.3333
$ws.Range("L116").Value = 3815.8333
$ws.Range("M116").Value = 43.66670000000022
$ws.Range("N116").Value = -10699.8333
$ws.Range("H122").Value = 938.5333000000001
$ws.Range("I122").Value = 950.0833
$ws.Range("J122").Value = 892.3333
$ws.Range("K122").Value = 2850.2499
$ws.Range("L122").Value = 2676.9999
$ws.Range("M122").Value = -400.2498999999998
$ws.Range("N122").Value = -7576.9999
$ws.Range("H125").Value = 2146701.5
$ws.Range("I125").Value = 2575196.8
$ws.Range("K125").Value = 23176771.2
$ws.Range("M125").Value = -23174311.2
$ws.Range("H132").Value = 1401.8043
$ws.Range("I132").Value = 1331.7
$ws.Range("K132").Value = 3995.1
$ws.Range("M132").Value = -1465.1
$ws.Range("H135").Value = 826.93335
$ws.Range("I135").Value = 826.93335
$ws.Range("K135").Value = 7442.40015
$ws.Range("M135").Value = -4907.40015
$ws.Range("H138").Value = 2701.14
$ws.Range("J138").Value = 3041.2163
$ws.Range("L138").Value = 9123.6489
$ws.Range("N138").Value = -19403.6489
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4579.3335
$ws.Range("I2").Value = 4017.182
$ws.Range("J2").Value = 5462.7144
$ws.Range("K2").Value = 4017.182
$ws.Range("L2").Value = 5462.7144
$ws.Range("M2").Value = -3904.182
$ws.Range("N2").Value = -5688.7144
$ws.Range("H110").Value = 7583.1665
$ws.Range("I110").Value = 6333.1113
$ws.Range("K110").Value = 6333.1113
$ws.Range("M110").Value = -4288.1113
$ws.Range("H116").Value = 4579.3335
$ws.Range("I116").Value = 4017.182
$ws.Range("J116").Value = 5462.7144
$ws.Range("K116").Value = 4017.182
$ws.Range("L116").Value = 5462.7144
$ws.Range("M116").Value = -1723.182
$ws.Range("N116").Value = -10050.7144
$ws.Range("H122").Value = 1030.2727
$ws.Range("I122").Value = 1030.2727
$ws.Range("K122").Value = 3090.8181
$ws.Range("M122").Value = -640.8181
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4579.3335
$ws.Range("I3").Value = 4017.182
$ws.Range("J3").Value = 5462.7144
$ws.Range("K3").Value = 4017.182
$ws.Range("L3").Value = 5462.7144
$ws.Range("M3").Value = -3903.182
$ws.Range("N3").Value = -5690.7144
$ws.Range("H20").Value = 3064.4736
$ws.Range("I20").Value = 2853.5
$ws.Range("J20").Value = 3298.889
$ws.Range("K20").Value = 2853.5
$ws.Range("L20").Value = 3298.889
$ws.Range("M20").Value = -2606.5
$ws.Range("N20").Value = -3792.889
$ws.Range("H134").Value = 3634.6
$ws.Range("I134").Value = 3322.0754
$ws.Range("K134").Value = 9966.226200000001
$ws.Range("M134").Value = -7431.226200000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 9850
$ws.Range("I6").Value = 17660.834
$ws.Range("K6").Value = 17660.834
$ws.Range("M6").Value = -17547.834
$ws.Range("H58").Value = 11682.477
$ws.Range("I58").Value = 7486.5
$ws.Range("J58").Value = 13360.866
$ws.Range("K58").Value = 7486.5
$ws.Range("L58").Value = 13360.866
$ws.Range("M58").Value = -7283.5
$ws.Range("N58").Value = -13766.866
$ws.Range("H99").Value = 4979.4
$ws.Range("I99").Value = 5249.25
$ws.Range("J99").Value = 3900
$ws.Range("K99").Value = 5249.25
$ws.Range("L99").Value = 3900
$ws.Range("M99").Value = -3751.25
$ws.Range("N99").Value = -6896
$ws.Range("H126").Value = 4979.4
$ws.Range("I126").Value = 5249.25
$ws.Range("J126").Value = 3900
$ws.Range("K126").Value = 15747.75
$ws.Range("L126").Value = 11700
$ws.Range("M126").Value = -13277.75
$ws.Range("N126").Value = -16640
$ws.Range("H132").Value = 3563.4707
$ws.Range("J132").Value = 5766.3335
$ws.Range("L132").Value = 17299.0005
$ws.Range("N132").Value = -22359.0005
$ws.Range("H134").Value = 3629.25
$ws.Range("I134").Value = 2221.2
$ws.Range("K134").Value = 6663.599999999999
$ws.Range("M134").Value = -4128.599999999999
$ws.Range("H136").Value = 11682.477
$ws.Range("I136").Value = 7486.5
$ws.Range("J136").Value = 13360.866
$ws.Range("K136").Value = 22459.5
$ws.Range("L136").Value = 40082.598
$ws.Range("M136").Value = -19909.5
$ws.Range("N136").Value = -45182.598
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 35992.855
$ws.Range("I17").Value = 412.5
$ws.Range("J17").Value = 83433.336
$ws.Range("K17").Value = 1237.5
$ws.Range("L17").Value = 250300.008
$ws.Range("M17").Value = -1068.5
$ws.Range("N17").Value = -250638.008
$ws.Range("H39").Value = 680
$ws.Range("J39").Value = 425
$ws.Range("L39").Value = 1275
$ws.Range("N39").Value = -1863
$ws.Range("H40").Value = 118.63158
$ws.Range("J40").Value = 171.25
$ws.Range("L40").Value = 685
$ws.Range("N40").Value = -823
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 610.0278
$ws.Range("I97").Value = 547.96295
$ws.Range("K97").Value = 547.96295
$ws.Range("M97").Value = -51.96294999999998
$ws.Range("H102").Value = 4614
$ws.Range("I102").Value = 3275.9167
$ws.Range("J102").Value = 8628.25
$ws.Range("K102").Value = 3275.9167
$ws.Range("L102").Value = 8628.25
$ws.Range("M102").Value = -1653.9167
$ws.Range("N102").Value = -11872.25
$ws.Range("H107").Value = 1267
$ws.Range("J107").Value = 1499.6
$ws.Range("L107").Value = 1499.6
$ws.Range("N107").Value = -5339.6
$ws.Range("H122").Value = 1953
$ws.Range("I122").Value = 1824.5294
$ws.Range("J122").Value = 2499
$ws.Range("K122").Value = 5473.5882
$ws.Range("L122").Value = 7497
$ws.Range("M122").Value = -3023.5882
$ws.Range("N122").Value = -12397
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6777.5
$ws.Range("I7").Value = 6799.8
$ws.Range("K7").Value = 6799.8
$ws.Range("M7").Value = -6687.8
$ws.Range("H40").Value = 4477.143
$ws.Range("I40").Value = 4556.6665
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 4556.6665
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -4420.6665
$ws.Range("N40").Value = -4272
$ws.Range("H61").Value = 2587.6875
$ws.Range("I61").Value = 2686.9333
$ws.Range("K61").Value = 2686.9333
$ws.Range("M61").Value = -2484.9333
$ws.Range("H68").Value = 4292.846
$ws.Range("J68").Value = 6358.143
$ws.Range("L68").Value = 6358.143
$ws.Range("N68").Value = -7856.143
$ws.Range("H71").Value = 4292.846
$ws.Range("J71").Value = 6358.143
$ws.Range("L71").Value = 31790.715
$ws.Range("N71").Value = -39278.715
$ws.Range("H113").Value = 2587.6875
$ws.Range("I113").Value = 2686.9333
$ws.Range("K113").Value = 2686.9333
$ws.Range("M113").Value = -516.9333000000001
$ws.Range("H126").Value = 6777.5
$ws.Range("I126").Value = 6799.8
$ws.Range("K126").Value = 20399.4
$ws.Range("M126").Value = -17929.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1886
$ws.Range("I107").Value = 1248.3334
$ws.Range("J107").Value = 3799
$ws.Range("K107").Value = 3745.0002
$ws.Range("L107").Value = 11397
$ws.Range("M107").Value = -1825.0002
$ws.Range("N107").Value = -15237
$ws.Range("H113").Value = 709.0476
$ws.Range("I113").Value = 573.8182
$ws.Range("J113").Value = 857.8
$ws.Range("K113").Value = 1721.4546
$ws.Range("L113").Value = 2573.4
$ws.Range("M113").Value = 448.5454
$ws.Range("N113").Value = -6913.4
$ws.Range("H126").Value = 5525.125
$ws.Range("I126").Value = 3530.5386
$ws.Range("K126").Value = 10591.6158
$ws.Range("M126").Value = -8121.6158
$ws.Range("H132").Value = 5144.048
$ws.Range("I132").Value = 4731.1353
$ws.Range("K132").Value = 14193.4059
$ws.Range("M132").Value = -11663.4059
